# Rename the only worksheet from "Sheet1" to "TigerSnus".
# Excel automatically updates any defined-name formulas (such as the
# hidden _FilterDatabase name) that reference the sheet by its old name.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "TigerSnus"

# Remove every defined name in the workbook except the built-in
# _FilterDatabase name, which must be kept (now pointing at the
# renamed sheet).
for ($i = $wb.Names.Count; $i -ge 1; $i--) {
    $n = $wb.Names.Item($i)
    if ($n.Name -notmatch "_FilterDatabase$") {
        $n.Delete()
    }
}
